# GMS Data Release 1
# The "referral_participant" data-dictionary sheet: rename the row that
# describes the participant id column from "patient_id" to
# "participant_id" (cell B4), matching the rest of the sheet's naming
# convention (see B5 "referral_id" etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

[void]$ws.Activate()
$ws.Cells.Item(4, 2).Value = "participant_id"

# Leave the selection where the edit was made, as Excel would after a
# manual edit of that cell.
[void]$ws.Range("B4").Select()
